# Update "想去人数" (number of people interested) values for a handful of
# events across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 16516
$wsExpo.Range("F6").Value = 15637
$wsExpo.Range("F30").Value = 5823

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 83

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 16516
$wsAll.Range("F6").Value = 15637
$wsAll.Range("F22").Value = 83
$wsAll.Range("F32").Value = 5823
